$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relabel the summary rows at the bottom of the table (B40:B43) ---
# Row 40: "Общая сумма" -> "Общая сумма, руб."
$ws.Range("B40").Value = "Общая сумма, руб."

# Row 41: "Средняя площадь" -> "Средняя площадь, кв.м."
$ws.Range("B41").Value = "Средняя площадь, кв.м."

# Row 42 ("Максимальный срок просрочки") keeps its text unchanged.

# Row 43: "Максимальная сумма" -> "Максимальная сумма, руб."
$ws.Range("B43").Value = "Максимальная сумма, руб."

# C43's formula switches from referencing the "Итого" column (K) to the
# "Сумма" column (E): MAX(K3:K38) -> MAX(E3:E38)
$ws.Range("C43").Formula = "=MAX(E3:E38)"

# --- Cosmetic view changes ---
# Zoom 70% -> 85%
$excel.ActiveWindow.Zoom = 85

# Selection moves from D41 to N14
[void]$ws.Range("N14").Select()

# Column B widens from ~29.375 to ~35.25 characters. The host's ColumnWidth
# setter snaps to the nearest 1/7th-character (pixel) increment, so feed it
# the pre-image of that rounding to land as close as possible to 35.25.
$ws.Columns("B").ColumnWidth = 34.571428571428573
